# Insert a new data row at row 47 (pushing the existing rows 47-132 down
# to 48-133, growing the sheet's used range from A1:R132 to A1:R133), then
# populate the newly inserted row with the new "Apio" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 47..132 down to 48..133, leaving a blank row 47 (inheriting
# the D-column's date number format from the row that used to sit there).
$ws.Rows(47).Insert()

# Populate the newly-inserted row 47 with the new record.
$ws.Cells.Item(47, 1).Value = 5
$ws.Cells.Item(47, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(47, 3).Value = 'Maule'
$ws.Cells.Item(47, 4).Value = 44495
$ws.Cells.Item(47, 5).Value = 7
$ws.Cells.Item(47, 6).Value = 100112017
$ws.Cells.Item(47, 7).Value = 'Apio'
$ws.Cells.Item(47, 8).Value = 'Americana (o)'
$ws.Cells.Item(47, 9).Value = 'Primera'
$ws.Cells.Item(47, 10).Value = 500
$ws.Cells.Item(47, 11).Value = 7000
$ws.Cells.Item(47, 12).Value = 7000
$ws.Cells.Item(47, 13).Value = 7000
$ws.Cells.Item(47, 14).Value = '$/docena de matas'
$ws.Cells.Item(47, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(47, 16).Value = 1167
$ws.Cells.Item(47, 17).Value = 6
$ws.Cells.Item(47, 18).Value = 'Hortaliza'
